$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.2881169905109251;  C = 0.3048912486333797; D = 3.223369029078222;  E = 13.86384647080068; G = 17.68022373902321 }
    3 = @{ B = 0.01253208636536152; C = 0.3048912486333797; D = 0.7210945179870265; E = 13.86384647080068; G = 14.90236432378645 }
    4 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    5 = @{ B = 0.04172184405617529; C = 0.04103571897497393; D = 0.7210945179870265; E = 0.5333859586016987; G = 1.337238039619874 }
    6 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 189.6080260415259;  E = 13.86384647080068; G = 208.3711874500482 }
    7 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 18.71679738969934;  E = 0.5333859586016987; G = 24.14949828602258 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
